$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 ("Left Bank Melbourne Restaurant and Cocktail Bar, 1 Southbank Blvd")
# was removed; every row below it shifts up by one.
$ws.Rows.Item(12).Delete()
